# Move the "Comment" column (column J, with its header, sub-header,
# data-type marker and format note) to the end of the table (after the
# "Energy" column, i.e. after column Q), shifting BombNumber..Energy
# (columns K:Q) one position to the left (J:P).
#
# "maj template comment a la fin" -> comment column goes to the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("J").Cut()
$ws.Columns("R").Insert()
